$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$targetRows = @(3,5,8,9,15,21,25,27,28,30,33,35,40,45,47,51,54,55,57,60,63,65,66,67,69,73,77,83,84,88,89,91,97,103,104,108,110,116,120,123,124,125,129,135,137,144,150,152,154,157,158,161,168,169,170,172,174,176,178,183,186,187,188,192)

$cols = @("G","H","I","J","K","L","M","N")

foreach ($r in $targetRows) {
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = "N/A"
    }
}
